$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.969.58"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.555.74"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'206.99"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D8").Value = "'22.12"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "1.777.29"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "1.555.32"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").Value = "26.966.93"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "'217.88"
$ws.Range("E18").Value = "  +1.88%  "
$ws.Range("E19").Value = "  +2.17%  "
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "'4.07"
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("D23").Value = "'9.22"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("D25").Value = "'154.45"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  +2.05%  "
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").Value = "1.427.17"
$ws.Range("E33").Value = "  +4.76%  "
$ws.Range("D34").Value = "'3.08"
$ws.Range("E34").Value = "  +4.67%  "
$ws.Range("E35").Value = "  +3.64%  "
$ws.Range("D36").Value = "'0.980"
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("D38").Value = "'0.0166"
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").Value = "'0.811"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("D41").Value = "'5.77"
$ws.Range("E41").Value = "  +3.11%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  +4.53%  "
$ws.Range("D44").Value = "'0.991"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").Value = "'64.39"
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("D47").Value = "1.691.13"
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("D48").Value = "'87.84"
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("D49").Value = "'0.0522"
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("D50").Value = "0.0₇0999"
$ws.Range("E50").Value = "  +3.39%  "
$ws.Range("D51").Value = "'0.0957"
$ws.Range("E51").Value = "  +0.70%  "
